# Natmi following Dr Hou advice
# Re-run of the Ptn-Sdc3 ligand-receptor pair analysis: the sender/receiver
# cluster grid now covers all 3 clusters (ECs, FAPs, sCs) x 3 clusters
# instead of just 2, and the underlying per-edge statistics are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T
$rows = @(
    @(2,  "ECs",  "Ptn", "Sdc3", "ECs",  1, 0.3333333333333333, 0.3206776666666667, 0.962033,           0.03001977461414601, 0.03001977461414601, 3, 1, 26.12444933333333, 78.37334799999999, 0.7238861157526749, 0.7238861157526749, 8.37752745516489,  75.397747096484,   0.02173089804120491, 0.02173089804120491),
    @(3,  "ECs",  "Ptn", "Sdc3", "FAPs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033,           0.03001977461414601, 0.03001977461414601, 3, 1, 3.818542,          11.455626,          0.1058085282850919, 0.1058085282850919, 1.224521138628667, 11.020690247658,   0.003176348171372953, 0.003176348171372953),
    @(4,  "ECs",  "Ptn", "Sdc3", "sCs",  1, 0.3333333333333333, 0.3206776666666667, 0.962033,           0.03001977461414601, 0.03001977461414601, 3, 1, 6.146179,          18.438537,          0.1703053559622332, 0.1703053559622332, 1.970942340635667, 17.738481065721,   0.005112528401568149, 0.005112528401568149),
    @(5,  "FAPs", "Ptn", "Sdc3", "ECs",  3, 1,                  5.752274333333333,  17.256823,          0.5384908178993973,  0.5384908178993975,  3, 1, 26.12444933333333, 78.37334799999999, 0.7238861157526749, 0.7238861157526749, 150.2749993726004, 1352.474994353404, 0.3898060265376758,  0.3898060265376758),
    @(6,  "FAPs", "Ptn", "Sdc3", "FAPs", 3, 1,                  5.752274333333333,  17.256823,          0.5384908178993973,  0.5384908178993975,  3, 1, 3.818542,          11.455626,          0.1058085282850919, 0.1058085282850919, 21.96530113735534, 197.687710236198,  0.05697692093697068, 0.05697692093697069),
    @(7,  "FAPs", "Ptn", "Sdc3", "sCs",  3, 1,                  5.752274333333333,  17.256823,          0.5384908178993973,  0.5384908178993975,  3, 1, 6.146179,          18.438537,          0.1703053559622332, 0.1703053559622332, 35.35450770977234, 318.190569387951,  0.09170787042475098, 0.091707870424751),
    @(8,  "sCs",  "Ptn", "Sdc3", "ECs",  3, 1,                  4.609262333333334,  13.827787,          0.4314894074864565,  0.4314894074864565,  3, 1, 26.12444933333333, 78.37334799999999, 0.7238861157526749, 0.7238861157526749, 120.4144402912084, 1083.729962620876, 0.3123491911737942,  0.3123491911737942),
    @(9,  "sCs",  "Ptn", "Sdc3", "FAPs", 3, 1,                  4.609262333333334,  13.827787,          0.4314894074864565,  0.4314894074864565,  3, 1, 3.818542,          11.455626,          0.1058085282850919, 0.1058085282850919, 17.60066180885133, 158.405956279662,  0.04565525917674829, 0.04565525917674829),
    @(10, "sCs",  "Ptn", "Sdc3", "sCs",  3, 1,                  4.609262333333334,  13.827787,          0.4314894074864565,  0.4314894074864565,  3, 1, 6.146179,          18.438537,          0.1703053559622332, 0.1703053559622332, 28.32935135862434, 254.964162227619,  0.07348495713591409, 0.07348495713591409)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
